$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 99) four times,
# inserting new rows 100-103 directly below it so the new rows inherit the
# same cell styles (column styles + the "Keep" highlight style in column G).
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(99).Copy() | Out-Null
    $ws.Rows.Item(100).Insert() | Out-Null
}

# Row 100: SeaBird / D / Station depth imputed from another site visit
$ws.Range("A100").Value = "SeaBird"
$ws.Range("B100").Value = "D"
$ws.Range("C100").Value = "Station depth imputed from another site visit"
$ws.Range("E100").Value = "D"
$ws.Range("F100").Value = "Station depth imputed"
$ws.Range("G100").Value = "Keep"

# Row 101: SeaBird / D / Station depth imputed from maximum CTD sample depth
$ws.Range("A101").Value = "SeaBird"
$ws.Range("B101").Value = "D"
$ws.Range("C101").Value = "Station depth imputed from maximum CTD sample depth"
$ws.Range("E101").Value = "D"
$ws.Range("F101").Value = "Station depth imputed"
$ws.Range("G101").Value = "Keep"

# Row 102: NOAActd / D / Station depth imputed from another site visit
$ws.Range("A102").Value = "NOAActd"
$ws.Range("B102").Value = "D"
$ws.Range("C102").Value = "Station depth imputed from another site visit"
$ws.Range("E102").Value = "D"
$ws.Range("F102").Value = "Station depth imputed"
$ws.Range("G102").Value = "Keep"

# Row 103: NOAActd / D / Station depth imputed from maximum CTD sample depth
$ws.Range("A103").Value = "NOAActd"
$ws.Range("B103").Value = "D"
$ws.Range("C103").Value = "Station depth imputed from maximum CTD sample depth"
$ws.Range("E103").Value = "D"
$ws.Range("F103").Value = "Station depth imputed"
$ws.Range("G103").Value = "Keep"

# Move the visible selection to reflect the new bottom of the sheet,
# matching what Excel saved after this edit.
$ws.Range("A104").Select() | Out-Null
